$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 0.9999970461623656
$ws.Cells.Item(2, 4).Value = 1.025507051647289
$ws.Cells.Item(2, 5).Value = 1.003512603051345
$ws.Cells.Item(2, 6).Value = 0.9980539377174651
$ws.Cells.Item(2, 9).Value = 1.028047310659273
$ws.Cells.Item(2, 10).Value = 1.005312350954267
$ws.Cells.Item(2, 11).Value = 1.028332420894459
$ws.Cells.Item(2, 12).Value = 1.006403886417997
$ws.Cells.Item(2, 13).Value = 1.000962071010937
$ws.Cells.Item(2, 14).Value = 1.005586813710726

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.002190696077171
$ws.Cells.Item(3, 4).Value = 1.026019325364456
$ws.Cells.Item(3, 5).Value = 1.005418250165748
$ws.Cells.Item(3, 6).Value = 1.000925108801209
$ws.Cells.Item(3, 9).Value = 1.028114196249099
$ws.Cells.Item(3, 10).Value = 1.007127342732227
$ws.Cells.Item(3, 11).Value = 1.02865281339994
$ws.Cells.Item(3, 12).Value = 1.008108944823934
$ws.Cells.Item(3, 13).Value = 1.003628619941815
$ws.Cells.Item(3, 14).Value = 1.006228286426888

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.003602290855779
$ws.Cells.Item(4, 4).Value = 1.026344013162027
$ws.Cells.Item(4, 5).Value = 1.006644382404554
$ws.Cells.Item(4, 6).Value = 1.002773591153557
$ws.Cells.Item(4, 9).Value = 1.028152502439901
$ws.Cells.Item(4, 10).Value = 1.008294031562883
$ws.Cells.Item(4, 11).Value = 1.028852644648525
$ws.Cells.Item(4, 12).Value = 1.009204923891012
$ws.Cells.Item(4, 13).Value = 1.005344591826161
$ws.Cells.Item(4, 14).Value = 1.006639578476929

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.004193894425286
$ws.Cells.Item(5, 4).Value = 1.026478881439717
$ws.Cells.Item(5, 5).Value = 1.007158223888107
$ws.Cells.Item(5, 6).Value = 1.003548520839981
$ws.Cells.Item(5, 9).Value = 1.028167412619662
$ws.Cells.Item(5, 10).Value = 1.008782696506916
$ws.Cells.Item(5, 11).Value = 1.028934856412294
$ws.Cells.Item(5, 12).Value = 1.009663960881819
$ws.Cells.Item(5, 13).Value = 1.006063784403245
$ws.Cells.Item(5, 14).Value = 1.006811593515251

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.004293121253613
$ws.Cells.Item(6, 4).Value = 1.026501430639359
$ws.Cells.Item(6, 5).Value = 1.007244405997112
$ws.Cells.Item(6, 6).Value = 1.003678509513614
$ws.Cells.Item(6, 9).Value = 1.028169846008701
$ws.Cells.Item(6, 10).Value = 1.008864640454659
$ws.Cells.Item(6, 11).Value = 1.028948554567337
$ws.Cells.Item(6, 12).Value = 1.009740935846062
$ws.Cells.Item(6, 13).Value = 1.006184412705869
$ws.Cells.Item(6, 14).Value = 1.006840423675391

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.003610203023348
$ws.Cells.Item(7, 4).Value = 1.026345821692126
$ws.Cells.Item(7, 5).Value = 1.006651254710766
$ws.Cells.Item(7, 6).Value = 1.002783954247485
$ws.Cells.Item(7, 9).Value = 1.028152706365015
$ws.Cells.Item(7, 10).Value = 1.008300568190578
$ws.Cells.Item(7, 11).Value = 1.028853750235325
$ws.Cells.Item(7, 12).Value = 1.009211064246002
$ws.Cells.Item(7, 13).Value = 1.005354210280489
$ws.Cells.Item(7, 14).Value = 1.006641880437511

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.000740052902658
$ws.Cells.Item(8, 4).Value = 1.02568157912909
$ws.Cells.Item(8, 5).Value = 1.004158087951371
$ws.Cells.Item(8, 6).Value = 0.9990262421512628
$ws.Cells.Item(8, 9).Value = 1.028070942867339
$ws.Cells.Item(8, 10).Value = 1.005927361266509
$ws.Cells.Item(8, 11).Value = 1.028442245416833
$ws.Cells.Item(8, 12).Value = 1.006981654404119
$ws.Cells.Item(8, 13).Value = 1.001865241835927
$ws.Cells.Item(8, 14).Value = 1.005804394172819

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 0.9956202764081897
$ws.Cells.Item(9, 4).Value = 1.024459342914767
$ws.Cells.Item(9, 5).Value = 0.9997098096606482
$ws.Cells.Item(9, 6).Value = 0.9923299142558761
$ws.Cells.Item(9, 9).Value = 1.027888906352644
$ws.Cells.Item(9, 10).Value = 1.001684422044947
$ws.Cells.Item(9, 11).Value = 1.027660033183786
$ws.Cells.Item(9, 12).Value = 1.002995500289088
$ws.Cells.Item(9, 13).Value = 0.9956418554368396
$ws.Cells.Item(9, 14).Value = 1.004299033365596

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 0.9921622425845256
$ws.Cells.Item(10, 4).Value = 1.023610008354521
$ws.Cells.Item(10, 5).Value = 0.9967047873794241
$ws.Cells.Item(10, 6).Value = 0.9878110444644356
$ws.Cells.Item(10, 9).Value = 1.027742192411911
$ws.Cells.Item(10, 10).Value = 0.9988121622532348
$ws.Cells.Item(10, 11).Value = 1.027100466945094
$ws.Cells.Item(10, 12).Value = 1.000296919877218
$ws.Cells.Item(10, 13).Value = 0.9914381395819487
$ws.Cells.Item(10, 14).Value = 1.003274668619613

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 0.9906535521406715
$ws.Cells.Item(11, 4).Value = 1.02323410637176
$ws.Cells.Item(11, 5).Value = 0.9953936374089645
$ws.Cells.Item(11, 6).Value = 0.9858403664321082
$ws.Cells.Item(11, 9).Value = 1.027672681232816
$ws.Cells.Item(11, 10).Value = 0.9975575053047278
$ws.Cells.Item(11, 11).Value = 1.026849189472056
$ws.Cells.Item(11, 12).Value = 0.9991181050659224
$ws.Cells.Item(11, 13).Value = 0.9896039559150833
$ws.Cells.Item(11, 14).Value = 1.002825974291096

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 0.9900913939142384
$ws.Cells.Item(12, 4).Value = 1.023093261465666
$ws.Cells.Item(12, 5).Value = 0.9949050722402973
$ws.Cells.Item(12, 6).Value = 0.9851061820500319
$ws.Cells.Item(12, 9).Value = 1.027645965106242
$ws.Cells.Item(12, 10).Value = 0.9970897735031714
$ws.Cells.Item(12, 11).Value = 1.026754508480711
$ws.Cells.Item(12, 12).Value = 0.9986786442991347
$ws.Cells.Item(12, 13).Value = 0.9889204816658863
$ws.Cells.Item(12, 14).Value = 1.002658519692851

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 0.9902120595381855
$ws.Cells.Item(13, 4).Value = 1.023123528240741
$ws.Cells.Item(13, 5).Value = 0.9950099418895325
$ws.Cells.Item(13, 6).Value = 0.9852637674690576
$ws.Cells.Item(13, 9).Value = 1.027651736342749
$ws.Cells.Item(13, 10).Value = 0.9971901811635615
$ws.Cells.Item(13, 11).Value = 1.026774878692961
$ws.Cells.Item(13, 12).Value = 0.9987729831400927
$ws.Cells.Item(13, 13).Value = 0.9890671890455449
$ws.Cells.Item(13, 14).Value = 1.00269447529006

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 0.9906071202911284
$ws.Cells.Item(14, 4).Value = 1.023222488916466
$ws.Cells.Item(14, 5).Value = 0.9953532842849065
$ws.Cells.Item(14, 6).Value = 0.9857797236426981
$ws.Cells.Item(14, 9).Value = 1.027670491146611
$ws.Cells.Item(14, 10).Value = 0.9975188773544146
$ws.Cells.Item(14, 11).Value = 1.026841390534216
$ws.Cells.Item(14, 12).Value = 0.9990818119466078
$ws.Cells.Item(14, 13).Value = 0.9895475046040829
$ws.Cells.Item(14, 14).Value = 1.002812148652134

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 0.9908502947745117
$ws.Cells.Item(15, 4).Value = 1.023283300548361
$ws.Cells.Item(15, 5).Value = 0.9955646225320298
$ws.Cells.Item(15, 6).Value = 0.9860973288618624
$ws.Cells.Item(15, 9).Value = 1.027681927861248
$ws.Cells.Item(15, 10).Value = 0.9977211715879997
$ws.Cells.Item(15, 11).Value = 1.026882192518696
$ws.Cells.Item(15, 12).Value = 0.9992718785883367
$ws.Cells.Item(15, 13).Value = 0.9898431519607224
$ws.Cells.Item(15, 14).Value = 1.002884545935586

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 0.9922621254230287
$ws.Cells.Item(16, 4).Value = 1.023634784689578
$ws.Cells.Item(16, 5).Value = 0.9967915900936591
$ws.Cells.Item(16, 6).Value = 0.9879415300188076
$ws.Cells.Item(16, 9).Value = 1.027746679784027
$ws.Cells.Item(16, 10).Value = 0.9988951946573905
$ws.Cells.Item(16, 11).Value = 1.02711695452225
$ws.Cells.Item(16, 12).Value = 1.00037493267123
$ws.Cells.Item(16, 13).Value = 0.9915595675619321
$ws.Cells.Item(16, 14).Value = 1.003304337241511

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 0.9931446518951077
$ws.Cells.Item(17, 4).Value = 1.023853086440131
$ws.Cells.Item(17, 5).Value = 0.9975585341355973
$ws.Cells.Item(17, 6).Value = 0.989094545834366
$ws.Cells.Item(17, 9).Value = 1.027785696793538
$ws.Cells.Item(17, 10).Value = 0.9996286609749276
$ws.Cells.Item(17, 11).Value = 1.027261812832631
$ws.Cells.Item(17, 12).Value = 1.001064055618971
$ws.Cells.Item(17, 13).Value = 0.992632438877268
$ws.Cells.Item(17, 14).Value = 1.003566273402792

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 0.993658324721283
$ws.Cells.Item(18, 4).Value = 1.02397963304905
$ws.Cells.Item(18, 5).Value = 0.9980049224938602
$ws.Cells.Item(18, 6).Value = 0.9897657377869169
$ws.Cells.Item(18, 9).Value = 1.027807877329898
$ws.Cells.Item(18, 10).Value = 1.000055426526088
$ws.Cells.Item(18, 11).Value = 1.027345439387843
$ws.Cells.Item(18, 12).Value = 1.001465017693937
$ws.Cells.Item(18, 13).Value = 0.9932568864796915
$ws.Cells.Item(18, 14).Value = 1.003718561842482

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 0.9938332907999199
$ws.Cells.Item(19, 4).Value = 1.024022648879112
$ws.Cells.Item(19, 5).Value = 0.9981569685386646
$ws.Cells.Item(19, 6).Value = 0.9899943717553462
$ws.Cells.Item(19, 9).Value = 1.027815342305293
$ws.Cells.Item(19, 10).Value = 1.000200765487866
$ws.Cells.Item(19, 11).Value = 1.027373806722036
$ws.Cells.Item(19, 12).Value = 1.001601568606567
$ws.Cells.Item(19, 13).Value = 0.9934695821076397
$ws.Cells.Item(19, 14).Value = 1.003770404919007

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 0.9930500783226616
$ws.Cells.Item(20, 4).Value = 1.023829745906674
$ws.Cells.Item(20, 5).Value = 0.9974763476563392
$ws.Cells.Item(20, 6).Value = 0.9889709776493107
$ws.Cells.Item(20, 9).Value = 1.027781570346447
$ws.Cells.Item(20, 10).Value = 0.9995500762588492
$ws.Cells.Item(20, 11).Value = 1.027246360539083
$ws.Cells.Item(20, 12).Value = 1.00099022216119
$ws.Cells.Item(20, 13).Value = 0.9925174691326852
$ws.Cells.Item(20, 14).Value = 1.003538221405817

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 0.9904908339126361
$ws.Cells.Item(21, 4).Value = 1.023193381061208
$ws.Cells.Item(21, 5).Value = 0.9952522215828185
$ws.Cells.Item(21, 6).Value = 0.9856278484838727
$ws.Cells.Item(21, 9).Value = 1.027664993056495
$ws.Cells.Item(21, 10).Value = 0.9974221317641361
$ws.Cells.Item(21, 11).Value = 1.026821841561943
$ws.Cells.Item(21, 12).Value = 0.9989909140153368
$ws.Cells.Item(21, 13).Value = 0.9894061243963045
$ws.Cells.Item(21, 14).Value = 1.002777518719308

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 0.9888715021031147
$ws.Cells.Item(22, 4).Value = 1.022786227006252
$ws.Cells.Item(22, 5).Value = 0.9938448567504581
$ws.Cells.Item(22, 6).Value = 0.9835131937094171
$ws.Cells.Item(22, 9).Value = 1.027586509990163
$ws.Cells.Item(22, 10).Value = 0.9960743677744542
$ws.Cells.Item(22, 11).Value = 1.026547147649342
$ws.Cells.Item(22, 12).Value = 0.9977246088080334
$ws.Cells.Item(22, 13).Value = 0.9874372626241351
$ws.Cells.Item(22, 14).Value = 1.002294657611214

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 0.9897309308844535
$ws.Cells.Item(23, 4).Value = 1.023002733729557
$ws.Cells.Item(23, 5).Value = 0.9945917945827535
$ws.Cells.Item(23, 6).Value = 0.9846354452697469
$ws.Cells.Item(23, 9).Value = 1.027628606217047
$ws.Cells.Item(23, 10).Value = 0.9967897932096135
$ws.Cells.Item(23, 11).Value = 1.026693504457306
$ws.Cells.Item(23, 12).Value = 0.998396795041804
$ws.Cells.Item(23, 13).Value = 0.9884822189254497
$ws.Cells.Item(23, 14).Value = 1.002551071293827

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 0.9930928154140612
$ws.Cells.Item(24, 4).Value = 1.023840294915894
$ws.Cells.Item(24, 5).Value = 0.9975134871437007
$ws.Cells.Item(24, 6).Value = 0.9890268169374212
$ws.Cells.Item(24, 9).Value = 1.027783436694908
$ws.Cells.Item(24, 10).Value = 0.9995855885604165
$ws.Cells.Item(24, 11).Value = 1.027253345444725
$ws.Cells.Item(24, 12).Value = 1.001023587384582
$ws.Cells.Item(24, 13).Value = 0.9925694231479877
$ws.Cells.Item(24, 14).Value = 1.003550898422906

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 0.9969515404821367
$ws.Cells.Item(25, 4).Value = 1.024781422877931
$ws.Cells.Item(25, 5).Value = 1.000866570495859
$ws.Cells.Item(25, 6).Value = 0.9940703804778075
$ws.Cells.Item(25, 9).Value = 1.027940447854657
$ws.Cells.Item(25, 10).Value = 1.002788820867278
$ws.Cells.Item(25, 11).Value = 1.02786898787746
$ws.Cells.Item(25, 12).Value = 1.004033090616584
$ws.Cells.Item(25, 13).Value = 0.9972601003021522
$ws.Cells.Item(25, 14).Value = 1.004691803431707
